# Apply updated cryptocurrency price/volume data to the worksheet.
# Cells whose new text would otherwise be auto-parsed by Excel as a
# number (losing formatting, e.g. "247.40" -> 247.4) are written with
# a leading single-quote so they stay literal text, matching the
# original workbook where every Price/Volume cell is stored as a string.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.242.56"
$ws.Range("E2").Value = "  +2.06%  "
$ws.Range("D3").Value = "2.004.03"
$ws.Range("E3").Value = "  +2.94%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'247.40"
$ws.Range("E5").Value = "  +1.78%  "
$ws.Range("E6").Value = "  +3.36%  "
$ws.Range("E7").Value = "  +5.04%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +3.17%  "
$ws.Range("D10").Value = "'0.0804"
$ws.Range("E10").Value = "  +2.75%  "
$ws.Range("E11").Value = "  +1.44%  "
$ws.Range("D12").Value = "'15.01"
$ws.Range("E12").Value = "  +9.86%  "
$ws.Range("D13").Value = "'22.74"
$ws.Range("E13").Value = "  +8.78%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "'0.849"
$ws.Range("E14").Value = "  +3.15%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.298.59"
$ws.Range("E15").Value = "  +2.80%  "
$ws.Range("D16").Value = "'5.49"
$ws.Range("E16").Value = "  +4.23%  "
$ws.Range("D17").Value = "2.005.17"
$ws.Range("E17").Value = "  +2.36%  "
$ws.Range("D18").Value = "37.179.12"
$ws.Range("E18").Value = "  +2.16%  "
$ws.Range("D19").Value = "'70.45"
$ws.Range("E19").Value = "  +1.52%  "
$ws.Range("E20").Value = "  +2.74%  "
$ws.Range("E21").Value = "  +4.30%  "
$ws.Range("D22").Value = "'230.84"
$ws.Range("E22").Value = "  +0.89%  "
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("E24").Value = "  +2.41%  "
$ws.Range("E25").Value = "  +0.76%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "'9.44"
$ws.Range("E26").Value = "  +4.34%  "
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").Value = "'0.144"
$ws.Range("E27").Value = "  +6.85%  "
$ws.Range("D28").Value = "'163.64"
$ws.Range("E28").Value = "  +2.44%  "
$ws.Range("D29").Value = "'19.73"
$ws.Range("E29").Value = "  +2.64%  "
$ws.Range("E30").Value = "  +14.43%  "
$ws.Range("E31").Value = "  +1.84%  "
$ws.Range("E32").Value = "  +3.82%  "
$ws.Range("D33").Value = "'0.0665"
$ws.Range("E33").Value = "  +10.03%  "
$ws.Range("E34").Value = "  +4.92%  "
$ws.Range("E35").Value = "  +6.61%  "
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("E37").Value = "  +2.39%  "
$ws.Range("D38").Value = "'3.28"
$ws.Range("E38").Value = "  -2.90%  "
$ws.Range("D39").Value = "'5.44"
$ws.Range("E39").Value = "  +5.32%  "
$ws.Range("D40").Value = "'0.0982"
$ws.Range("E40").Value = "  +1.54%  "
$ws.Range("E41").Value = "  +0.87%  "
$ws.Range("D42").Value = "'0.0216"
$ws.Range("E42").Value = "  +3.39%  "
$ws.Range("E43").Value = "  +3.25%  "
$ws.Range("D44").Value = "'16.77"
$ws.Range("E44").Value = "  +7.27%  "
$ws.Range("D45").Value = "'91.29"
$ws.Range("E45").Value = "  +4.98%  "
$ws.Range("D46").Value = "1.379.95"
$ws.Range("E46").Value = "  +1.36%  "
$ws.Range("E47").Value = "  +3.58%  "
$ws.Range("D48").Value = "'7.29"
$ws.Range("E48").Value = "  +2.67%  "
$ws.Range("E49").Value = "  +17.93%  "
$ws.Range("D50").Value = "'2.85"
$ws.Range("E50").Value = "  +0.87%  "
$ws.Range("E51").Value = "  +6.42%  "
